$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- csv column (D) ---
$ws.Cells.Item(2, 4).Value = "compare(expected,actual,failFast)"
$ws.Cells.Item(3, 4).Value = "compareExtended(var,profile,expected,actual)"
$ws.Cells.Item(4, 4).Value = "convertExcel(excel,worksheet,csvFile)"
$ws.Cells.Item(5, 4).Value = "fromExcel(excel,worksheet,csvFile)"

# --- io column (I) ---
$ws.Cells.Item(2, 9).Value = "assertEqual(expected,actual)"
$ws.Cells.Item(3, 9).Value = "assertNotEqual(expected,actual)"
$ws.Cells.Item(4, 9).Value = "assertReadableFile(file,minByte)"
$ws.Cells.Item(5, 9).Value = "compare(expected,actual,failFast)"
$ws.Cells.Item(6, 9).Value = "copyFiles(source,target)"
$ws.Cells.Item(7, 9).Value = "count(var,path,pattern)"
$ws.Cells.Item(8, 9).Value = "deleteFiles(location,recursive)"
$ws.Cells.Item(9, 9).Value = "filter(source,target,matchPattern)"
$ws.Cells.Item(10, 9).Value = "makeDirectory(source)"
$ws.Cells.Item(11, 9).Value = "moveFiles(source,target)"
$ws.Cells.Item(12, 9).Value = "readFile(var,file)"
$ws.Cells.Item(13, 9).Value = "readProperty(var,file,property)"
$ws.Cells.Item(14, 9).Value = "saveDiff(var,expected,actual)"
$ws.Cells.Item(15, 9).Value = "saveFileMeta(var,file)"
$ws.Cells.Item(16, 9).Value = "saveMatches(var,path,filePattern)"
$ws.Cells.Item(17, 9).Value = "unzip(zipFile,target)"
$ws.Cells.Item(18, 9).Value = "validate(var,profile,inputFile)"
$ws.Cells.Item(19, 9).Value = "writeFile(file,content,append)"
$ws.Cells.Item(20, 9).Value = "writeFileAsIs(file,content,append)"
$ws.Cells.Item(21, 9).Value = "writeProperty(file,property,value)"
$ws.Cells.Item(22, 9).Value = "zip(filePattern,zipFile)"

# --- json column (K) ---
$ws.Cells.Item(2, 11).Value = "addOrReplace(json,jsonpath,input,var)"
$ws.Cells.Item(3, 11).Value = "assertCorrectness(json,schema)"
$ws.Cells.Item(4, 11).Value = "assertElementCount(json,jsonpath,count)"
$ws.Cells.Item(5, 11).Value = "assertElementNotPresent(json,jsonpath)"
$ws.Cells.Item(6, 11).Value = "assertElementPresent(json,jsonpath)"
$ws.Cells.Item(7, 11).Value = "assertValue(json,jsonpath,expected)"
$ws.Cells.Item(8, 11).Value = "assertValues(json,jsonpath,array,exactOrder)"
$ws.Cells.Item(9, 11).Value = "assertWellformed(json)"
$ws.Cells.Item(10, 11).Value = "fromCsv(csv,header,jsonFile)"
$ws.Cells.Item(11, 11).Value = "storeCount(json,jsonpath,var)"
$ws.Cells.Item(12, 11).Value = "storeValue(json,jsonpath,var)"
$ws.Cells.Item(13, 11).Value = "storeValues(json,jsonpath,var)"

# --- webalert column (V) ---
$ws.Cells.Item(2, 22).Value = "accept()"
$ws.Cells.Item(3, 22).Value = "assertPresent()"
$ws.Cells.Item(4, 22).Value = "assertText(text,matchBy)"
$ws.Cells.Item(5, 22).Value = "dismiss()"
$ws.Cells.Item(6, 22).Value = "replyCancel(text)"
$ws.Cells.Item(7, 22).Value = "replyOK(text)"
$ws.Cells.Item(8, 22).Value = "storeText(var)"

# --- web column (U) ---
$ws.Cells.Item(2, 21).Value = "assertAndClick(locator,label)"
$ws.Cells.Item(3, 21).Value = "assertAttribute(locator,attrName,value)"
$ws.Cells.Item(4, 21).Value = "assertAttributeContains(locator,attrName,contains)"
$ws.Cells.Item(5, 21).Value = "assertAttributeNotContains(locator,attrName,contains)"
$ws.Cells.Item(6, 21).Value = "assertAttributeNotPresent(locator,attrName)"
$ws.Cells.Item(7, 21).Value = "assertAttributePresent(locator,attrName)"
$ws.Cells.Item(8, 21).Value = "assertChecked(locator)"
$ws.Cells.Item(9, 21).Value = "assertContainCount(locator,text,count)"
$ws.Cells.Item(10, 21).Value = "assertCssNotPresent(locator,property)"
$ws.Cells.Item(11, 21).Value = "assertCssPresent(locator,property,value)"
$ws.Cells.Item(12, 21).Value = "assertElementByAttributes(nameValues)"
$ws.Cells.Item(13, 21).Value = "assertElementByText(locator,text)"
$ws.Cells.Item(14, 21).Value = "assertElementCount(locator,count)"
$ws.Cells.Item(15, 21).Value = "assertElementNotPresent(locator)"
$ws.Cells.Item(16, 21).Value = "assertElementPresent(locator)"
$ws.Cells.Item(17, 21).Value = "assertFocus(locator)"
$ws.Cells.Item(18, 21).Value = "assertFrameCount(count)"
$ws.Cells.Item(19, 21).Value = "assertFramePresent(frameName)"
$ws.Cells.Item(20, 21).Value = "assertIECompatMode()"
$ws.Cells.Item(21, 21).Value = "assertIENavtiveMode()"
$ws.Cells.Item(22, 21).Value = "assertLinkByLabel(label)"
$ws.Cells.Item(23, 21).Value = "assertNotChecked(locator)"
$ws.Cells.Item(24, 21).Value = "assertNotFocus(locator)"
$ws.Cells.Item(25, 21).Value = "assertNotText(locator,text)"
$ws.Cells.Item(26, 21).Value = "assertNotVisible(locator)"
$ws.Cells.Item(27, 21).Value = "assertOneMatch(locator)"
$ws.Cells.Item(28, 21).Value = "assertScrollbarHNotPresent(locator)"
$ws.Cells.Item(29, 21).Value = "assertScrollbarHPresent(locator)"
$ws.Cells.Item(30, 21).Value = "assertScrollbarVNotPresent(locator)"
$ws.Cells.Item(31, 21).Value = "assertScrollbarVPresent(locator)"
$ws.Cells.Item(32, 21).Value = "assertTable(locator,row,column,text)"
$ws.Cells.Item(33, 21).Value = "assertText(locator,text)"
$ws.Cells.Item(34, 21).Value = "assertTextContains(locator,text)"
$ws.Cells.Item(35, 21).Value = "assertTextCount(locator,text,count)"
$ws.Cells.Item(36, 21).Value = "assertTextList(locator,list,ignoreOrder)"
$ws.Cells.Item(37, 21).Value = "assertTextMatches(text,minMatch,scrollTo)"
$ws.Cells.Item(38, 21).Value = "assertTextNotPresent(text)"
$ws.Cells.Item(39, 21).Value = "assertTextOrder(locator,descending)"
$ws.Cells.Item(40, 21).Value = "assertTextPresent(text)"
$ws.Cells.Item(41, 21).Value = "assertTitle(text)"
$ws.Cells.Item(42, 21).Value = "assertValue(locator,value)"
$ws.Cells.Item(43, 21).Value = "assertValueOrder(locator,descending)"
$ws.Cells.Item(44, 21).Value = "assertVisible(locator)"
$ws.Cells.Item(45, 21).Value = "checkAll(locator)"
$ws.Cells.Item(46, 21).Value = "clearLocalStorage()"
$ws.Cells.Item(47, 21).Value = "click(locator)"
$ws.Cells.Item(48, 21).Value = "clickAndWait(locator,waitMs)"
$ws.Cells.Item(49, 21).Value = "clickByLabel(label)"
$ws.Cells.Item(50, 21).Value = "clickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(51, 21).Value = "close()"
$ws.Cells.Item(52, 21).Value = "closeAll()"
$ws.Cells.Item(53, 21).Value = "deselectMulti(locator,array)"
$ws.Cells.Item(54, 21).Value = "dismissInvalidCert()"
$ws.Cells.Item(55, 21).Value = "dismissInvalidCertPopup()"
$ws.Cells.Item(56, 21).Value = "doubleClick(locator)"
$ws.Cells.Item(57, 21).Value = "doubleClickAndWait(locator,waitMs)"
$ws.Cells.Item(58, 21).Value = "doubleClickByLabel(label)"
$ws.Cells.Item(59, 21).Value = "doubleClickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(60, 21).Value = "dragAndDrop(fromLocator,toLocator)"
$ws.Cells.Item(61, 21).Value = "editLocalStorage(key,value)"
$ws.Cells.Item(62, 21).Value = "executeScript(var,script)"
$ws.Cells.Item(63, 21).Value = "focus(locator)"
$ws.Cells.Item(64, 21).Value = "goBack()"
$ws.Cells.Item(65, 21).Value = "goBackAndWait()"
$ws.Cells.Item(66, 21).Value = "maximizeWindow()"
$ws.Cells.Item(67, 21).Value = "mouseOver(locator)"
$ws.Cells.Item(68, 21).Value = "open(url)"
$ws.Cells.Item(69, 21).Value = "openAndWait(url,waitMs)"
$ws.Cells.Item(70, 21).Value = "refresh()"
$ws.Cells.Item(71, 21).Value = "refreshAndWait()"
$ws.Cells.Item(72, 21).Value = "resizeWindow(width,height)"
$ws.Cells.Item(73, 21).Value = "saveAllWindowIds(var)"
$ws.Cells.Item(74, 21).Value = "saveAllWindowNames(var)"
$ws.Cells.Item(75, 21).Value = "saveAttribute(var,locator,attrName)"
$ws.Cells.Item(76, 21).Value = "saveCount(var,locator)"
$ws.Cells.Item(77, 21).Value = "saveElement(var,locator)"
$ws.Cells.Item(78, 21).Value = "saveElements(var,locator)"
$ws.Cells.Item(79, 21).Value = "saveLocalStorage(var,key)"
$ws.Cells.Item(80, 21).Value = "saveLocation(var)"
$ws.Cells.Item(81, 21).Value = "savePageAs(var,sessionIdName,url)"
$ws.Cells.Item(82, 21).Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Cells.Item(83, 21).Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Cells.Item(84, 21).Value = "saveText(var,locator)"
$ws.Cells.Item(85, 21).Value = "saveTextArray(var,locator)"
$ws.Cells.Item(86, 21).Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Cells.Item(87, 21).Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Cells.Item(88, 21).Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Cells.Item(89, 21).Value = "saveValue(var,locator)"
$ws.Cells.Item(90, 21).Value = "scrollLeft(locator,pixel)"
$ws.Cells.Item(91, 21).Value = "scrollRight(locator,pixel)"
$ws.Cells.Item(92, 21).Value = "scrollTo(locator)"
$ws.Cells.Item(93, 21).Value = "select(locator,text)"
$ws.Cells.Item(94, 21).Value = "selectFrame(locator)"
$ws.Cells.Item(95, 21).Value = "selectMulti(locator,array)"
$ws.Cells.Item(96, 21).Value = "selectMultiOptions(locator)"
$ws.Cells.Item(97, 21).Value = "selectText(locator)"
$ws.Cells.Item(98, 21).Value = "selectWindow(winId)"
$ws.Cells.Item(99, 21).Value = "selectWindowAndWait(winId,waitMs)"
$ws.Cells.Item(100, 21).Value = "selectWindowByIndex(index)"
$ws.Cells.Item(101, 21).Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Cells.Item(102, 21).Value = "toggleSelections(locator)"
$ws.Cells.Item(103, 21).Value = "type(locator,value)"
$ws.Cells.Item(104, 21).Value = "typeKeys(locator,value)"
$ws.Cells.Item(105, 21).Value = "uncheckAll(locator)"
$ws.Cells.Item(106, 21).Value = "unselectAllText()"
$ws.Cells.Item(107, 21).Value = "upload(fieldLocator,file)"
$ws.Cells.Item(108, 21).Value = "verifyContainText(locator,text)"
$ws.Cells.Item(109, 21).Value = "verifyText(locator,text)"
$ws.Cells.Item(110, 21).Value = "wait(waitMs)"
$ws.Cells.Item(111, 21).Value = "waitForElementPresent(locator)"
$ws.Cells.Item(112, 21).Value = "waitForPopUp(winId,waitMs)"
$ws.Cells.Item(113, 21).Value = "waitForTextPresent(text)"
$ws.Cells.Item(114, 21).Value = "waitForTitle(text)"
# --- update defined names ranges to reflect new extents ---
$wb.Names.Item("csv").RefersTo = "='#system'!`$D`$2:`$D`$5"
$wb.Names.Item("io").RefersTo = "='#system'!`$I`$2:`$I`$22"
$wb.Names.Item("json").RefersTo = "='#system'!`$K`$2:`$K`$13"
$wb.Names.Item("web").RefersTo = "='#system'!`$U`$2:`$U`$114"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$V`$2:`$V`$8"
